$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 201) holds the "Förändrad" date, stored as the
# serial number 45189 (2023-09-20). Bump it by one day to 45190 (2023-09-21).
$ws.Range("C2:C201").Value = 45190
